# Weekly fruit/vegetable price refresh for "Bruselas (repollito)" — Vega Central
# Mapocho de Santiago. Every existing data row (2-38) gets new date / volume /
# price figures, and a brand-new row (39) is appended, growing the used range
# from A1:R38 to A1:R39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D=44726; J=28; K=24000; L=24000; M=24000; P=1600 },
    @{ Row=3; D=44442; J=28; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=4; D=44343; J=26; K=23000; L=24000; M=23500; P=1567 },
    @{ Row=5; D=44349; J=21; K=24000; L=25000; M=24524; P=1635 },
    @{ Row=6; D=44400; J=16; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=7; D=44425; J=25; K=24000; L=25000; M=24520; P=1635 },
    @{ Row=8; D=44685; J=20; K=25000; L=25000; M=25000; P=1667 },
    @{ Row=9; D=44421; J=18; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=10; D=44351; J=34; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=11; D=44705; J=35; K=26000; L=26000; M=26000; P=1733 },
    @{ Row=12; D=44411; J=34; K=25000; L=26000; M=25500; P=1700 },
    @{ Row=13; D=44719; J=43; K=17000; L=18000; M=17512; P=1167 },
    @{ Row=14; D=44460; J=25; K=24000; L=25000; M=24480; P=1632 },
    @{ Row=15; D=44432; J=34; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=16; D=44707; J=30; K=26000; L=26000; M=26000; P=1733 },
    @{ Row=17; D=44463; J=25; K=24000; L=25000; M=24480; P=1632 },
    @{ Row=18; D=44390; J=34; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=19; D=44341; J=36; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=20; D=44677; J=34; K=25000; L=26000; M=25500; P=1700 },
    @{ Row=21; D=44336; J=34; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=22; D=44446; J=34; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=23; D=44455; J=18; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=24; D=44428; J=16; K=25000; L=26000; M=25500; P=1700 },
    @{ Row=25; D=44714; J=52; K=18000; L=20000; M=19000; P=1267 },
    @{ Row=26; D=44708; J=25; K=26000; L=26000; M=26000; P=1733 },
    @{ Row=27; D=44413; J=25; K=24000; L=25000; M=24480; P=1632 },
    @{ Row=28; D=44727; J=28; K=24000; L=24000; M=24000; P=1600 },
    @{ Row=29; D=44453; J=25; K=25000; L=26000; M=25520; P=1701 },
    @{ Row=30; D=44329; J=25; K=23000; L=23000; M=23000; P=1533 },
    @{ Row=31; D=44406; J=25; K=24000; L=25000; M=24520; P=1635 },
    @{ Row=32; D=44435; J=34; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=33; D=44680; J=36; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=34; D=44418; J=16; K=25000; L=26000; M=25500; P=1700 },
    @{ Row=35; D=44706; J=30; K=26000; L=26000; M=26000; P=1733 },
    @{ Row=36; D=44383; J=25; K=13000; L=14000; M=13480; P=899 },
    @{ Row=37; D=44449; J=18; K=24000; L=25000; M=24500; P=1633 },
    @{ Row=38; D=44385; J=25; K=14000; L=15000; M=14480; P=965 },
    @{ Row=39; D=44397; J=34; K=23000; L=24000; M=23500; P=1567 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio $/Kg
}

# Row 39 is brand new, so the columns that are constant across every other
# row in the sheet need to be written explicitly too.
$ws.Range("D39").NumberFormat = $ws.Range("D38").NumberFormat
$ws.Range("A39").Value = 9
$ws.Range("B39").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C39").Value = 'Metropolitana'
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = 100112035
$ws.Range("G39").Value = 'Bruselas (repollito)'
$ws.Range("H39").Value = 'Sin especificar'
$ws.Range("I39").Value = 'Primera'
$ws.Range("N39").Value = '$/malla 15 kilos'
$ws.Range("O39").Value = 'Hijuelas'
$ws.Range("Q39").Value = 15
$ws.Range("R39").Value = 'Hortaliza'
